# "fix(gui) step 1 and 2"
# Daily price-list refresh on the "CABLE ACERO Y ALAMBRE ROPA" sheet:
#   - bump the date in A1 by one day (45308 -> 45309, i.e. 2024-01-17 -> 2024-01-18)
#   - update the two price cells (D22, D37) to their new values

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A1").Value = 45309
$ws.Range("D22").Value = 141
$ws.Range("D37").Value = 196
